# Insert a new weekly record as row 49 ("Hortaliza, Terminal Hortofrutícola
# Agro Chillán - Apio" daily-logic sheet). All existing rows from 49 onward
# shift down by one (old row 49 -> new row 50, ..., old row 120 -> new row
# 121); the sheet's used range grows from A1:R120 to A1:R121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 49..120 down to 50..121, inserting a blank row at 49.
$ws.Rows("49:49").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A49").Value2 = 7
$ws.Range("B49").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value2 = "Ñuble"
$ws.Range("D49").Value2 = 44467
$ws.Range("E49").Value2 = 16
$ws.Range("F49").Value2 = 100112017
$ws.Range("G49").Value2 = "Apio"
$ws.Range("H49").Value2 = "Americana (o)"
$ws.Range("I49").Value2 = "Primera"
$ws.Range("J49").Value2 = 160
$ws.Range("K49").Value2 = 8000
$ws.Range("L49").Value2 = 9000
$ws.Range("M49").Value2 = 8500
$ws.Range("N49").Value2 = "$/docena de matas"
$ws.Range("O49").Value2 = "Provincia del Elquí"
$ws.Range("P49").Value2 = 1417
$ws.Range("Q49").Value2 = 6
$ws.Range("R49").Value2 = "Hortaliza"
